$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.424.39'
$ws.Range("E2").Value = '  +2.05%  '
$ws.Range("D3").Value = '3.388.00'
$ws.Range("E3").Value = '  +1.55%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.99'
$ws.Range("E5").Value = '  +0.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.79'
$ws.Range("E6").Value = '  +2.13%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.595'
$ws.Range("E8").Value = '  +0.33%  '
$ws.Range("E9").Value = '  +9.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.590'
$ws.Range("E10").Value = '  +1.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '48.23'
$ws.Range("E11").Value = '  +0.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000285'
$ws.Range("E12").Value = '  +4.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '683.25'
$ws.Range("E13").Value = '  -2.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.64'
$ws.Range("E14").Value = '  +2.90%  '
$ws.Range("D15").Value = '3.944.83'
$ws.Range("E15").Value = '  +1.78%  '
$ws.Range("D16").Value = '69.506.03'
$ws.Range("E16").Value = '  +2.17%  '
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.121'
$ws.Range("E17").Value = '  +1.27%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.391.05'
$ws.Range("E18").Value = '  +1.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.68'
$ws.Range("E19").Value = '  +1.18%  '
$ws.Range("E20").Value = '  +1.54%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.909'
$ws.Range("E21").Value = '  +1.63%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.20'
$ws.Range("E22").Value = '  +1.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.34'
$ws.Range("E23").Value = '  -1.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '102.23'
$ws.Range("E24").Value = '  +1.79%  '
$ws.Range("E26").Value = '  +0.86%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.65'
$ws.Range("E27").Value = '  +2.38%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '33.86'
$ws.Range("E28").Value = '  +1.88%  '
$ws.Range("E29").Value = '  +3.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.93'
$ws.Range("E30").Value = '  -1.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.09'
$ws.Range("E31").Value = '  +0.88%  '
$ws.Range("B32").Value = 'dogwifhat'
$ws.Range("C32").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.64'
$ws.Range("E32").Value = '  +8.70%  '
$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '552.54'
$ws.Range("E33").Value = '  -3.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.106'
$ws.Range("E34").Value = '  +0.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '58.29'
$ws.Range("E35").Value = '  +1.89%  '
$ws.Range("E36").Value = '  +0.16%  '
$ws.Range("D37").Value = '3.659.47'
$ws.Range("E37").Value = '  -2.21%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.139'
$ws.Range("E38").Value = '  +4.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.55'
$ws.Range("E39").Value = '  +0.42%  '
$ws.Range("D40").Value = '0.0₃0725'
$ws.Range("E40").Value = '  +7.10%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.28'
$ws.Range("E41").Value = '  +4.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.70'
$ws.Range("E42").Value = '  +2.64%  '
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0427'
$ws.Range("E43").Value = '  +4.57%  '
$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.338'
$ws.Range("E44").Value = '  +1.07%  '
$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.33'
$ws.Range("E45").Value = '  +0.47%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.66'
$ws.Range("E46").Value = '  +1.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.129'
$ws.Range("E47").Value = '  +0.37%  '
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("E49").Value = '  +4.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '129.36'
$ws.Range("E50").Value = '  -0.75%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.66'
$ws.Range("E51").Value = '  -4.15%  '
